$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043008604067959
$ws.Range("D2").Value = 1.036728855291699
$ws.Range("E2").Value = 1.050379973544527
$ws.Range("F2").Value = 1.059593659591669
$ws.Range("I2").Value = 1.039711868729749
$ws.Range("J2").Value = 1.048081359011294
$ws.Range("K2").Value = 1.039521765596109
$ws.Range("L2").Value = 1.05313441241268
$ws.Range("M2").Value = 1.062322739536633
$ws.Range("N2").Value = 1.019850971105894
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044103394530057
$ws.Range("D3").Value = 1.037250246682439
$ws.Range("E3").Value = 1.051359533167427
$ws.Range("F3").Value = 1.060671129212168
$ws.Range("I3").Value = 1.039932449134138
$ws.Range("J3").Value = 1.048822099949598
$ws.Range("K3").Value = 1.039853583086355
$ws.Range("L3").Value = 1.053925976040393
$ws.Range("M3").Value = 1.063213803428208
$ws.Range("N3").Value = 1.020102494061603
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044811976759894
$ws.Range("D4").Value = 1.037587837719835
$ws.Range("E4").Value = 1.051993887745394
$ws.Range("F4").Value = 1.061368976404196
$ws.Range("I4").Value = 1.040074095910541
$ws.Range("J4").Value = 1.049301023533619
$ws.Range("K4").Value = 1.040067793559723
$ws.Range("L4").Value = 1.054438067228384
$ws.Range("M4").Value = 1.063790434252316
$ws.Range("N4").Value = 1.020264979662287
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045109908061339
$ws.Range("D5").Value = 1.037729811360415
$ws.Range("E5").Value = 1.052260693164548
$ws.Range("F5").Value = 1.061662507294455
$ws.Range("I5").Value = 1.04013338451634
$ws.Range("J5").Value = 1.049502270632579
$ws.Range("K5").Value = 1.040157727762183
$ws.Range("L5").Value = 1.054653325411766
$ws.Range("M5").Value = 1.064032862300714
$ws.Range("N5").Value = 1.020333224641996
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045159934588733
$ws.Range("D6").Value = 1.037753652274269
$ws.Range("E6").Value = 1.052305498158192
$ws.Range("F6").Value = 1.061711801559738
$ws.Range("I6").Value = 1.040143324096326
$ws.Range("J6").Value = 1.049536055485459
$ws.Range("K6").Value = 1.040172821058902
$ws.Range("L6").Value = 1.054689466739099
$ws.Range("M6").Value = 1.064073567760654
$ws.Range("N6").Value = 1.020344679518999
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.044815957561826
$ws.Range("D7").Value = 1.037589734582338
$ws.Range("E7").Value = 1.051997452330492
$ws.Range("F7").Value = 1.0613728979659
$ws.Range("I7").Value = 1.040074889148616
$ws.Range("J7").Value = 1.049303712970333
$ws.Range("K7").Value = 1.040068995736929
$ws.Range("L7").Value = 1.054440943617053
$ws.Range("M7").Value = 1.06379367353911
$ws.Range("N7").Value = 1.020265891806578
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043378556645825
$ws.Range("D8").Value = 1.036905016443673
$ws.Range("E8").Value = 1.050710913698658
$ws.Range("F8").Value = 1.059957660292417
$ws.Range("I8").Value = 1.039786639069587
$ws.Range("J8").Value = 1.048331775589814
$ws.Range("K8").Value = 1.039634007460016
$ws.Range("L8").Value = 1.053401946340562
$ws.Range("M8").Value = 1.062623867667314
$ws.Range("N8").Value = 1.019936029571964
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040847032691273
$ws.Range("D9").Value = 1.035700167930533
$ws.Range("E9").Value = 1.04844781768299
$ws.Range("F9").Value = 1.057468836342651
$ws.Range("I9").Value = 1.039270420835465
$ws.Range("J9").Value = 1.046616151749506
$ws.Range("K9").Value = 1.038863724109309
$ws.Range("L9").Value = 1.051570316696272
$ws.Range("M9").Value = 1.060562934571793
$ws.Range("N9").Value = 1.019352735012046
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.039160240285105
$ws.Range("D10").Value = 1.034898165252285
$ws.Range("E10").Value = 1.046941757174747
$ws.Range("F10").Value = 1.055812996406627
$ws.Range("I10").Value = 1.03892072262494
$ws.Range("J10").Value = 1.045470426569999
$ws.Range("K10").Value = 1.038347705177004
$ws.Range("L10").Value = 1.050348712311289
$ws.Range("M10").Value = 1.059189267635125
$ws.Range("N10").Value = 1.018962511594702
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038430045435851
$ws.Range("D11").Value = 1.034551197302421
$ws.Range("E11").Value = 1.04629025199474
$ws.Range("F11").Value = 1.055096801512137
$ws.Range("I11").Value = 1.038767985842994
$ws.Range("J11").Value = 1.044973846509021
$ws.Range("K11").Value = 1.03812368013332
$ws.Range("L11").Value = 1.049819622428933
$ws.Range("M11").Value = 1.058594524396468
$ws.Range("N11").Value = 1.018793219310591
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.038158847566999
$ws.Range("D12").Value = 1.034422365017488
$ws.Range("E12").Value = 1.046048348550119
$ws.Range("F12").Value = 1.054830894242078
$ws.Range("I12").Value = 1.038711055248649
$ws.Range("J12").Value = 1.044789323351561
$ws.Range("K12").Value = 1.038040380003266
$ws.Range("L12").Value = 1.049623075791813
$ws.Range("M12").Value = 1.058373619839432
$ws.Range("N12").Value = 1.018730288178165
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.038217019128038
$ws.Range("D13").Value = 1.034449997829352
$ws.Range("E13").Value = 1.046100233378356
$ws.Range("F13").Value = 1.054887926850861
$ws.Range("I13").Value = 1.038723275981272
$ws.Range("J13").Value = 1.044828907418923
$ws.Range("K13").Value = 1.038058252101125
$ws.Range("L13").Value = 1.049665236576815
$ws.Range("M13").Value = 1.058421004183947
$ws.Range("N13").Value = 1.018743789310409
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038407627547335
$ws.Range("D14").Value = 1.034540547017975
$ws.Range("E14").Value = 1.046270254250443
$ws.Range("F14").Value = 1.055074819098442
$ws.Range("I14").Value = 1.038763283966892
$ws.Range("J14").Value = 1.044958595220405
$ws.Range("K14").Value = 1.038116796294493
$ws.Range("L14").Value = 1.049803376212156
$ws.Range("M14").Value = 1.058576264161753
$ws.Range("N14").Value = 1.018788018391248
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038525071510203
$ws.Range("D15").Value = 1.034596343636686
$ws.Range("E15").Value = 1.046375022257715
$ws.Range("F15").Value = 1.055189985402522
$ws.Range("I15").Value = 1.038787908063082
$ws.Range("J15").Value = 1.045038490716568
$ws.Range("K15").Value = 1.038152855763158
$ws.Range("L15").Value = 1.049888486070198
$ws.Range("M15").Value = 1.058671926225678
$ws.Range("N15").Value = 1.018815262971331
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.039208705016029
$ws.Range("D16").Value = 1.034921198900092
$ws.Range("E16").Value = 1.046985008687264
$ws.Range("F16").Value = 1.055860544647272
$ws.Range("I16").Value = 1.03893083157335
$ws.Range("J16").Value = 1.045503372970716
$ws.Range("K16").Value = 1.038362560706739
$ws.Range("L16").Value = 1.050383823586051
$ws.Range("M16").Value = 1.0592287401258
$ws.Range("N16").Value = 1.018973740168471
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039637582643728
$ws.Range("D17").Value = 1.035125054476712
$ws.Range("E17").Value = 1.047367805514918
$ws.Range("F17").Value = 1.056281381518331
$ws.Range("I17").Value = 1.039020131777112
$ws.Range("J17").Value = 1.045794854504768
$ws.Range("K17").Value = 1.038493946834652
$ws.Range("L17").Value = 1.050694501689764
$ws.Range("M17").Value = 1.059578031808623
$ws.Range("N17").Value = 1.019073062317384
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039887758796704
$ws.Range("D18").Value = 1.03524398927799
$ws.Range("E18").Value = 1.047591145107153
$ws.Range("F18").Value = 1.05652692522429
$ws.Range("I18").Value = 1.039072092166613
$ws.Range("J18").Value = 1.0459648250714
$ws.Range("K18").Value = 1.038570525606724
$ws.Range("L18").Value = 1.05087570270578
$ws.Range("M18").Value = 1.059781773874597
$ws.Range("N18").Value = 1.019130964079106
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039973065657229
$ws.Range("D19").Value = 1.035284547900148
$ws.Range("E19").Value = 1.047667308425874
$ws.Range("F19").Value = 1.056610662288602
$ws.Range("I19").Value = 1.039089787783681
$ws.Range("J19").Value = 1.046022772896246
$ws.Range("K19").Value = 1.038596627378582
$ws.Range("L19").Value = 1.050937485489089
$ws.Range("M19").Value = 1.05985124569469
$ws.Range("N19").Value = 1.019150701794511
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039591566163209
$ws.Range("D20").Value = 1.035103179663259
$ws.Range("E20").Value = 1.047326728760964
$ws.Range("F20").Value = 1.056236221773409
$ws.Range("I20").Value = 1.039010563835406
$ws.Range("J20").Value = 1.045763586020131
$ws.Range("K20").Value = 1.038479856188981
$ws.Range("L20").Value = 1.050661170138048
$ws.Range("M20").Value = 1.059540555481801
$ws.Range("N20").Value = 1.01906240921394
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03835149732299
$ws.Range("D21").Value = 1.034513881229201
$ws.Range("E21").Value = 1.046220184729437
$ws.Range("F21").Value = 1.055019780693566
$ws.Range("I21").Value = 1.038751508058328
$ws.Range("J21").Value = 1.044920407361662
$ws.Range("K21").Value = 1.03809955889174
$ws.Range("L21").Value = 1.049762698059741
$ws.Range("M21").Value = 1.058530543704006
$ws.Range("N21").Value = 1.018774995365855
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03757198435155
$ws.Range("D22").Value = 1.034143638836303
$ws.Range("E22").Value = 1.04552500387944
$ws.Range("F22").Value = 1.054255646734738
$ws.Range("I22").Value = 1.038587487702396
$ws.Range("J22").Value = 1.044389855423486
$ws.Range("K22").Value = 1.037859946407079
$ws.Range("L22").Value = 1.049197682515814
$ws.Range("M22").Value = 1.057895564543484
$ws.Range("N22").Value = 1.018594006651591
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037985203261759
$ws.Range("D23").Value = 1.034339885027687
$ws.Range("E23").Value = 1.045893480491987
$ws.Range("F23").Value = 1.054660663119889
$ws.Range("I23").Value = 1.038674546138938
$ws.Range("J23").Value = 1.04467115012651
$ws.Range("K23").Value = 1.037987017104066
$ws.Range("L23").Value = 1.049497218418256
$ws.Range("M23").Value = 1.058232173763461
$ws.Range("N23").Value = 1.018689978727774
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039612358970757
$ws.Range("D24").Value = 1.035113063859505
$ws.Range("E24").Value = 1.047345289390348
$ws.Range("F24").Value = 1.056256627281199
$ws.Range("I24").Value = 1.039014887568531
$ws.Range("J24").Value = 1.045777715043418
$ws.Range("K24").Value = 1.038486223319973
$ws.Range("L24").Value = 1.05067623126931
$ws.Range("M24").Value = 1.059557489401654
$ws.Range("N24").Value = 1.019067222989103
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.041501333641924
$ws.Range("D25").Value = 1.036011439130826
$ws.Range("E25").Value = 1.049032412698586
$ws.Range("F25").Value = 1.058111662716698
$ws.Range("I25").Value = 1.039404855610615
$ws.Range("J25").Value = 1.047060030263275
$ws.Range("K25").Value = 1.039063304069577
$ws.Range("L25").Value = 1.052043928814872
$ws.Range("M25").Value = 1.061095685063177
$ws.Range("N25").Value = 1.019503770993307

Write-Host "Updated vm_pu values for Case_1_7 (380 kV case)"
